# edit.ps1
#
# Applies the commit "Modificacion de un archivo":
#   - Paragraph 1 ("UNIVERSIDAD") becomes centered and gets
#     " TECNICA PARTICULAR DE LOJA" appended (as its own run).
#   - A new paragraph containing "Version 1.0" is inserted right after it;
#     the "_GoBack" bookmark (originally trailing paragraph 1) ends up at
#     the end of this new paragraph.

$d = $word.ActiveDocument

# --- 1. Append " TECNICA PARTICULAR DE LOJA" after "UNIVERSIDAD" -----------
$p1 = $d.Paragraphs(1)
$insertStart = $p1.Range.End - 1   # just after "UNIVERSIDAD", before the ¶
$appendRange = $d.Range($insertStart, $insertStart)
$appendText = " TECNICA PARTICULAR DE LOJA"
$appendRange.InsertAfter($appendText)

# Force the appended text into its own <w:r> (same formatting would
# otherwise be silently merged back into the first run): bracket it with a
# throw-away bookmark, then remove the bookmark, leaving the run split
# intact but with no extra markup.
$newTextRange = $d.Range($insertStart, $insertStart + $appendText.Length)
$d.Bookmarks.Add("TEMP_SPLIT_MARK", $newTextRange)
$d.Bookmarks("TEMP_SPLIT_MARK").Delete()

# --- 2. Split off a new paragraph for "Version 1.0" -------------------------
# Do this *before* centering paragraph 1, so the new paragraph does not
# inherit the centered alignment.
$p1 = $d.Paragraphs(1)
$p1End = $p1.Range.End
$splitRange = $d.Range($p1End, $p1End)
$splitRange.InsertParagraphAfter()

# The "_GoBack" bookmark used to trail paragraph 1; remove it there and
# recreate it at the end of the new (now second) paragraph, after the
# "Version 1.0" text we are about to type.
$d.Bookmarks("_GoBack").Delete()

$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter("Version 1.0")

# Work around collapsed-range bookmarks misbehaving at the very last valid
# document position: pad with a throw-away character, anchor the bookmark
# right before it, then delete the padding.
$endPos = $d.Content.End
$padRange = $d.Range($endPos - 1, $endPos - 1)
$padRange.InsertAfter("X")

$bmPos = $d.Content.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padPos = $d.Content.End - 2
$d.Range($padPos, $padPos + 1).Delete()

# --- 3. Center paragraph 1 ("UNIVERSIDAD TECNICA PARTICULAR DE LOJA") ------
$p1 = $d.Paragraphs(1)
$p1.Alignment = 1   # wdAlignParagraphCenter

Write-Output $d.Content.Text
